# Add a new column R (2021 data) to the 9.2.1 table, mirroring the
# existing column Q formatting, then leave the selection on T3 as
# recorded in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- R2 (blank separator cell, same look as Q2) -----------------------
# Q2 uses the "thin bottom border, no extra fill" look; recreate the
# bottom border directly so we don't have to rely on copy/paste for a
# blank cell.
$r2 = $ws.Range("R2")
$r2.BorderAround(1, -4138, 1) | Out-Null
$r2.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft   -> none
$r2.Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> none
$r2.Borders.Item(10).LineStyle = -4142  # xlEdgeRight  -> none

# --- R3 (header, year 2021) --------------------------------------------
$ws.Range("Q3").Copy() | Out-Null
$ws.Range("R3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R3").Value = 2021

# --- R4 (GVA share percentage, 2021 value) ------------------------------
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R4").Value = 13.5

# --- R5 (GVA per capita, 2021 value) ------------------------------------
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R5").Value = 15.1

$excel.CutCopyMode = 0

# --- final selection, matches the saved workbook ------------------------
$ws.Range("T3").Select() | Out-Null
